$d = $word.ActiveDocument
$d.TrackRevisions = $false

$old = "> As discussed previously, our book extends the I-O techniques found in the literature. We explicitly incorporate the First and Second Laws of Thermodynamics into the revisions of the I-O accounting equations to allow embodied energy to accumulate and depreciate in each sector of the economy.  This enhancement to the I-O methodology is necessary to explore the effect on economic development as significant transitions in energy inputs and technology occur. "

$new = "As discussed previously, our book extends the I-O techniques found in the literature by We explicitly incorporate first-principles of Thermodynamics to derive I-O accounting equations that allow accumulation and depreciation of embodied energy in sectors of the economy.  This enhancement to the I-O methodology is necessary to address questions of economic development and energy transitions. "

$d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
